$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# operation / [kW_el*h/a] and [kW_el*h] for conversion_1 (rows 4-5)
$ws.Range("D4").Value = 18732.19023536615
$ws.Range("D5").Value = 18732.19023536615

# operation / [kW_el*h/a] and [kW_el*h] for conversion_2 (rows 9-10)
$ws.Range("D9").Value = 1485.548159853576
$ws.Range("D10").Value = 1485.548159853576

# operation / [kW_el*h/a] and [kW_el*h] for conversion_3 (rows 14-15)
$ws.Range("D14").Value = 1485.40976463414
$ws.Range("D15").Value = 1485.40976463414
